$d = $word.ActiveDocument

# --- Remove the old "_GoBack" bookmark first (it currently sits on the
#     "Nappy input cancel /YA" paragraph) so that, once we re-create it
#     further down on the "Nappy record" paragraph, there is never a
#     name clash between two "_GoBack" bookmarks. -----------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Locate the "Nappy record /Tim" paragraph ---------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Nappy record*") {
        $target = $cand
        break
    }
}

$full = $target.Range

# Replace the paragraph (runs + pPr + paragraph mark) with the fully
# color-formatted version, add the trailing " Done" run and move the
# "_GoBack" bookmark to the end of this paragraph.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="003B2252" w:rsidRDefault="003B2252" w:rsidP="00AC769B">
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:rPr><w:color w:val="76923C" w:themeColor="accent3" w:themeShade="BF"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:color w:val="76923C" w:themeColor="accent3" w:themeShade="BF"/></w:rPr><w:t>N</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="76923C" w:themeColor="accent3" w:themeShade="BF"/></w:rPr><w:t>appy record</w:t></w:r>
<w:r w:rsidR="006B6F63"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="76923C" w:themeColor="accent3" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> /Tim</w:t></w:r>
<w:r><w:rPr><w:color w:val="76923C" w:themeColor="accent3" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> Done</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$full.InsertXML($xml)
